# Regenerate experiment task-order sheets with newly randomized stim orders.
$wb = $excel.ActiveWorkbook

# --- Position 1: "GNG_TO-16512555641582" -> "vSAT_TO-16515889981628168" ---
$ws = $wb.Worksheets.Item(1)
$ws.Range("B2").Value = 'SAT_stims-165158899810035.csv'
$ws.Range("B3").Value = 'vSAT_stims-16515889981471925.csv'
$ws.Range("B4").Value = 'SAT_stims-16515889981159394.csv'
$ws.Range("B5").Value = 'vSAT_stims-16515889981315637.csv'
$ws.Name = 'vSAT_TO-16515889981628168'

# --- Position 2: "NB_TO-16512555667555287" -> "TOL_TO-1651588998209692" ---
$ws = $wb.Worksheets.Item(2)
$ws.Range("A8:B10").EntireRow.Delete()
$ws.Range("B2").Value = 'MM_stims-16515889981784422.csv'
$ws.Range("B3").Value = 'ZM_stims-16515889981628168.csv'
$ws.Range("B4").Value = 'MM_stims-16515889981940675.csv'
$ws.Range("B5").Value = 'ZM_stims-16515889981784422.csv'
$ws.Range("B6").Value = 'MM_stims-1651588998209692.csv'
$ws.Range("B7").Value = 'ZM_stims-16515889981940675.csv'
$ws.Name = 'TOL_TO-1651588998209692'

# --- Position 3: "RS_TO-16512555667622964" -> "NB_TO-1651588999266219" ---
$ws = $wb.Worksheets.Item(3)
$ws.Range("A2").Copy($ws.Range("A4:A10"))
$ws.Range("B2").Value = 'TB-16515889992506258.csv'
$ws.Range("B3").Value = 'OB-16515889989941664.csv'
$ws.Range("A4").Value = 2
$ws.Range("B4").Value = 'TB-16515889990879178.csv'
$ws.Range("A5").Value = 3
$ws.Range("B5").Value = 'OB-16515889985701027.csv'
$ws.Range("A6").Value = 4
$ws.Range("B6").Value = 'ZB-match_9-16515889984585156.csv'
$ws.Range("A7").Value = 5
$ws.Range("B7").Value = 'TB-16515889990254169.csv'
$ws.Range("A8").Value = 6
$ws.Range("B8").Value = 'OB-16515889989160776.csv'
$ws.Range("A9").Value = 7
$ws.Range("B9").Value = 'ZB-match_9-16515889984428573.csv'
$ws.Range("A10").Value = 8
$ws.Range("B10").Value = 'ZB-match_8-16515889983781765.csv'
$ws.Name = 'NB_TO-1651588999266219'

# --- Position 4: "TOL_TO-16512555668194845" -> "RS_TO-1651588999266219" ---
$ws = $wb.Worksheets.Item(4)
$ws.Range("A4:B7").EntireRow.Delete()
$ws.Range("B2").Value = 'eyes open'
$ws.Range("B3").Value = 'eyes closed'
$ws.Name = 'RS_TO-1651588999266219'

# --- Position 5: "vSAT_TO-16512555669073684" -> "GNG_TO-16515889993012643" ---
$ws = $wb.Worksheets.Item(5)
$ws.Range("B2").Value = 'go_stims-1651588999266219.csv'
$ws.Range("B3").Value = 'GNG_stims-16515889992818458.csv'
$ws.Range("B4").Value = 'go_stims-16515889992834368.csv'
$ws.Range("B5").Value = 'GNG_stims-16515889993012643.csv'
$ws.Name = 'GNG_TO-16515889993012643'
